# Online Enrollment Check List - mark items complete/incomplete, add new checklist rows
# "for company wide testing using pretend payment"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Row 15 ("CONVERGE") gets UN-checked: A15:D15 lose their X/date/checkmarks
#    and revert to the blank "unchecked" look. Copy the unchecked formatting
#    that currently lives on row 8 (A8:D8) BEFORE row 8 itself is modified.
# ---------------------------------------------------------------------------
$ws.Range("A8:D8").Copy()
$ws.Range("A15:D15").PasteSpecial($xlPasteFormats)
$ws.Range("A15:D15").ClearContents()

# E15 keeps its text ("CONVERGE") but switches to the plain unchecked style
# (same formatting as D50, a blank unchecked cell with style index 4).
$ws.Range("D50").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Row 8 ("Move CC authorize to CC Payment Screen") gets CHECKED off:
#    copy the checked formatting from row 3 (A3:D3) into A8:D8, then fill in
#    the X mark and completion date, keeping D8's existing text.
# ---------------------------------------------------------------------------
$ws.Range("A3:D3").Copy()
$ws.Range("A8:D8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A8").Value2 = "X"
$ws.Range("B8").Value2 = 45889

# ---------------------------------------------------------------------------
# 3) Row 52 ("Kids names need to be upper case") gets CHECKED off the same way.
# ---------------------------------------------------------------------------
$ws.Range("A3:D3").Copy()
$ws.Range("A52:D52").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A52").Value2 = "X"
$ws.Range("B52").Value2 = 45889
$ws.Range("D52").Value2 = "Kids names need to be upper case"

# ---------------------------------------------------------------------------
# 4) New note row 54: "add (applicable taxes not included)" in column E,
#    wrapped text, taller row.
# ---------------------------------------------------------------------------
$ws.Range("E54").Value2 = "add (applicable taxes not included)`n"
$ws.Range("E54").WrapText = $true
$ws.Rows(54).RowHeight = 30

# ---------------------------------------------------------------------------
# 5) New checklist row 55: "Both contracts need to say Enrollment Fee and
#    show the value" in column D (same blank/unchecked style as neighbors).
# ---------------------------------------------------------------------------
$ws.Range("D55").Value2 = "Both contracts need to say Enrollment Fee and show the value"

# ---------------------------------------------------------------------------
# 6) Update the view: scrolled down further, new active selection.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("H49").Select()
